$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.082.90'
$ws.Range("E2").Value = '  +2.75%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.300.27'
$ws.Range("E3").Value = '  +1.67%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.26'
$ws.Range("E5").Value = '  +1.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.69'
$ws.Range("E6").Value = '  +5.34%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.537'
$ws.Range("E7").Value = '  +2.26%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.521'
$ws.Range("E9").Value = '  +6.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.28'
$ws.Range("E10").Value = '  +4.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0826'
$ws.Range("E11").Value = '  +4.91%  '
$ws.Range("E12").Value = '  +0.83%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.15'
$ws.Range("E13").Value = '  +8.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.657.19'
$ws.Range("E14").Value = '  +1.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.01'
$ws.Range("E15").Value = '  +4.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.304.57'
$ws.Range("E16").Value = '  +1.95%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.809'
$ws.Range("E17").Value = '  +2.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.997.13'
$ws.Range("E18").Value = '  +2.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.49'
$ws.Range("E19").Value = '  +1.60%  '
$ws.Range("D20").Value = '0.0₃0924'
$ws.Range("E20").Value = '  +2.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.08'
$ws.Range("E21").Value = '  +2.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.42'
$ws.Range("E22").Value = '  +0.63%  '
$ws.Range("E23").Value = '  +1.21%  '
$ws.Range("E24").Value = '  +4.93%  '
$ws.Range("E25").Value = '  +1.85%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.35'
$ws.Range("E27").Value = '  +3.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '38.46'
$ws.Range("E28").Value = '  +5.21%  '
$ws.Range("E29").Value = '  +1.81%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.25'
$ws.Range("E31").Value = '  +4.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.33'
$ws.Range("E32").Value = '  +2.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.11'
$ws.Range("E34").Value = '  -0.64%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.78'
$ws.Range("E35").Value = '  +3.95%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0739'
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.40'
$ws.Range("E37").Value = '  +0.34%  '
$ws.Range("E38").Value = '  +1.78%  '
$ws.Range("E39").Value = '  +1.82%  '
$ws.Range("E40").Value = '  +0.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.22'
$ws.Range("E41").Value = '  +5.46%  '
$ws.Range("E42").Value = '  -3.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.969.96'
$ws.Range("E43").Value = '  -0.20%  '
$ws.Range("E44").Value = '  +2.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.10'
$ws.Range("E45").Value = '  +1.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.03'
$ws.Range("E46").Value = '  +3.56%  '
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.64'
$ws.Range("E48").Value = '  +4.86%  '
$ws.Range("E49").Value = '  +16.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.526.18'
$ws.Range("E50").Value = '  +1.49%  '
$ws.Range("E51").Value = '  +2.25%  '
